$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 new values (previously held by row 7)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = "2 2 1 99"
$ws.Range("J5").Value = "OUTRAS"
$ws.Range("K5").Value = 2600000000
$ws.Range("L5").Value = 2600000000
$ws.Range("M5").Value = 2600000000

# Row 7 new values (previously held by row 5)
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = "2 1 2 4"
$ws.Range("J7").Value = "KFW"
$ws.Range("K7").Value = 350000000
$ws.Range("L7").Value = 350000000
$ws.Range("M7").Value = 350000000
